$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 49 (pushes existing rows 49-73 down to 50-74)
$ws.Rows.Item(49).EntireRow.Insert()

# Populate the new row 49 with the new weekly price entry
$ws.Range("A49").Value = 10
$ws.Range("B49").Value = "Vega Modelo de Temuco"
$ws.Range("C49").Value = "La Araucanía"
$ws.Range("D49").Value = 44777
$ws.Range("E49").Value = 9
$ws.Range("F49").Value = 300000001
$ws.Range("G49").Value = "Rabanito"
$ws.Range("H49").Value = "Sin especificar"
$ws.Range("I49").Value = "Primera"
$ws.Range("J49").Value = 55
$ws.Range("K49").Value = 10000
$ws.Range("L49").Value = 10000
$ws.Range("M49").Value = 10000
$ws.Range("N49").Value = "`$/docena de paquetes"
$ws.Range("O49").Value = "Provincia de Cautín"
$ws.Range("P49").Value = 833
$ws.Range("Q49").Value = 12
$ws.Range("R49").Value = "Hortaliza"
